# Update view-count ("想去人数") and price ("最低票价") figures on the
# "展览" and "全部类型" sheets, which contain duplicated data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> hashtable of column letter -> new value
$updates = @{
    2  = @{ F = 42 }
    4  = @{ F = 60 }
    5  = @{ F = 497 }
    6  = @{ F = 1466; G = 60 }
    7  = @{ F = 815 }
    8  = @{ F = 104 }
    9  = @{ F = 198 }
    10 = @{ G = 50 }
    11 = @{ F = 196 }
    12 = @{ F = 113 }
    13 = @{ F = 169 }
    14 = @{ F = 150 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $addr = "$col$row"
            $ws.Range($addr).Value = $cols[$col]
        }
    }
}
